$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.312.53"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").Value = "1.860.06"
$ws.Range("E3").Value = "  +0.04%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.6995"
$ws.Range("E5").Value = "  -0.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.48"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07771"
$ws.Range("E8").Value = "  -2.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3037"
$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.66"
$ws.Range("E10").Value = "  +5.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08146"
$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("D12").Value = "1.851.87"
$ws.Range("E12").Value = "  -1.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.202"
$ws.Range("E13").Value = "  +0.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7124"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.12"
$ws.Range("E15").Value = "  -0.13%  "

$ws.Range("D16").Value = "29.272.09"
$ws.Range("E16").Value = "  +0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.772"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "241.64"
$ws.Range("E18").Value = "  +2.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007753"
$ws.Range("E19").Value = "  -0.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.13"
$ws.Range("E20").Value = "  -2.01%  "

$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").Value = "2.084.87"
$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("E23").Value = "  +0.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.500"
$ws.Range("E24").Value = "  +0.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.89"
$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.868"
$ws.Range("E26").Value = "  -1.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1427"
$ws.Range("E27").Value = "  -1.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.02"
$ws.Range("E28").Value = "  -0.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.893"
$ws.Range("E29").Value = "  -5.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.369"
$ws.Range("E30").Value = "  -4.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.475"
$ws.Range("E31").Value = "  -0.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.281"
$ws.Range("E32").Value = "  -2.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.019"
$ws.Range("E33").Value = "  -0.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05148"
$ws.Range("E34").Value = "  -1.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.178"
$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7029"
$ws.Range("E36").Value = "  -0.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9964"
$ws.Range("E37").Value = "  -0.29%  "

$ws.Range("E38").Value = "  +0.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01841"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("E40").Value = "  -0.97%  "

$ws.Range("D41").Value = "1.174.93"
$ws.Range("E41").Value = "  +2.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9136"
$ws.Range("E42").Value = "  -1.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.989"
$ws.Range("E43").Value = "  +0.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.50"
$ws.Range("E44").Value = "  +0.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4228"
$ws.Range("E45").Value = "  -1.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.13"
$ws.Range("E47").Value = "  -1.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5353"
$ws.Range("E48").Value = "  -1.19%  "

$ws.Range("E49").Value = "  -2.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.120"
$ws.Range("E50").Value = "  -0.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.914"
$ws.Range("E51").Value = "  -0.47%  "
